$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename D1 header from "client code" to "clientNo"
$ws.Range("D1").Value = "clientNo"

# Add new column E header "groupNo"
$ws.Range("E1").Value = "groupNo"

# Update selection to D2
$ws.Range("D2").Select()
